$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Apply updated coin price/listing data (symbol list update).
# Numeric-looking values are written with a leading apostrophe to keep them
# as text (matching the source inlineStr cells), then the quote-prefix style
# that Excel applies is reset back to Normal so no stray formatting is left
# behind.

# Row 2
$ws.Range("D2").Value = "'245.06"
$ws.Range("D2").Style = "Normal"

# Row 3
$ws.Range("D3").Value = "'22.00"
$ws.Range("D3").Style = "Normal"

# Row 4
$ws.Range("D4").Value = "'5.337"
$ws.Range("D4").Style = "Normal"

# Row 5
$ws.Range("D5").Value = "'0.05956"
$ws.Range("D5").Style = "Normal"

# Row 6
$ws.Range("D6").Value = "'3.395"
$ws.Range("D6").Style = "Normal"

# Row 7
$ws.Range("D7").Value = "'6.388"
$ws.Range("D7").Style = "Normal"

# Row 8
$ws.Range("D8").Value = "'0.8069"
$ws.Range("D8").Style = "Normal"

# Row 9
$ws.Range("D9").Value = "'0.9664"
$ws.Range("D9").Style = "Normal"

# Row 10
$ws.Range("D10").Value = "'0.1426"
$ws.Range("D10").Style = "Normal"

# Row 11
$ws.Range("D11").Value = "'0.07396"
$ws.Range("D11").Style = "Normal"

# Row 12
$ws.Range("D12").Value = "'0.03414"
$ws.Range("D12").Style = "Normal"

# Row 13
$ws.Range("D13").Value = "'0.03052"
$ws.Range("D13").Style = "Normal"

# Row 14
$ws.Range("D14").Value = "'0.09404"
$ws.Range("D14").Style = "Normal"

# Row 15
$ws.Range("D15").Value = "'3.992"
$ws.Range("D15").Style = "Normal"

# Row 16
$ws.Range("D16").Value = "'0.001598"
$ws.Range("D16").Style = "Normal"

# Row 17
$ws.Range("D17").Value = "'0.04816"
$ws.Range("D17").Style = "Normal"

# Row 18
$ws.Range("B18").Value = "TigerCash"
$ws.Range("C18").Value = "https://coinranking.com/coin/6hIn06L2+tigercash-tch"
$ws.Range("D18").Value = "'0.006197"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "17TigerCashTCH"

# Row 19
$ws.Range("B19").Value = "HotbitToken"
$ws.Range("C19").Value = "https://coinranking.com/coin/uQJB8Ocu8lTb+hotbittoken-htb"
$ws.Range("D19").Value = "'0.005133"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "18HotbitTokenHTB"

# Row 20
$ws.Range("B20").Value = "BitKan"
$ws.Range("C20").Value = "https://coinranking.com/coin/RDOsLDgvY-AXe+bitkan-kan"
$ws.Range("D20").Value = "'0.0009853"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "19BitKanKAN"

# Row 21
$ws.Range("B21").Value = "NitroEx"
$ws.Range("C21").Value = "https://coinranking.com/coin/8oiZw6gwYhC+nitroex-ntx"
$ws.Range("D21").Value = "'0.00009704"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "20NitroExNTX"

# Row 22
$ws.Range("B22").Value = "LEO"
$ws.Range("C22").Value = "https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"
$ws.Range("D22").Value = "'3.745"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "21LEOLEO"

# Row 23
$ws.Range("B23").Value = "BTSEToken"
$ws.Range("C23").Value = "https://coinranking.com/coin/EOSL_JJKNMEr+btsetoken-btse"
$ws.Range("D23").Value = "'2.167"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "22BTSETokenBTSE"

# Row 24
$ws.Range("B24").Value = "One"
$ws.Range("C24").Value = "https://coinranking.com/coin/6Lga5NiXX3rT+one-one"
$ws.Range("D24").Value = "'0.0005913"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "23OneONE"

# Row 27
$ws.Range("D27").Value = "'0.0002462"
$ws.Range("D27").Style = "Normal"

# Row 40
$ws.Range("D40").Value = "'0.03908"
$ws.Range("D40").Style = "Normal"

# Row 41
$ws.Range("B41").Value = "KickToken"
$ws.Range("C41").Value = "https://coinranking.com/coin/F_Yv9Cu7pPL3Y+kicktoken-kick"
$ws.Range("D41").Value = "'0.006523"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "40KickTokenKICK"

# Row 42
$ws.Range("B42").Value = "BKEXToken"
$ws.Range("C42").Value = "https://coinranking.com/coin/IPeThtYgk+bkextoken-bkk"
$ws.Range("D42").Value = "'0.1069"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "41BKEXTokenBKK"

# Row 43
$ws.Range("B43").Value = "CEJI"
$ws.Range("C43").Value = "https://coinranking.com/coin/SbKjCVJCh+ceji-ceji"
$ws.Range("D43").Value = "'0.003001"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "42CEJICEJI"

# Row 44
$ws.Range("D44").Value = "'0.005843"
$ws.Range("D44").Style = "Normal"

# Row 45
$ws.Range("D45").Value = "'0.00005313"
$ws.Range("D45").Style = "Normal"

# Row 46
$ws.Range("D46").Value = "'0.00000000750"
$ws.Range("D46").Style = "Normal"

# Row 47
$ws.Range("D47").Value = "'0.8504"
$ws.Range("D47").Style = "Normal"

# Row 48
$ws.Range("D48").Value = "'0.03648"
$ws.Range("D48").Style = "Normal"

# Row 50
$ws.Range("D50").Value = "'0.01010"
$ws.Range("D50").Style = "Normal"
